# Update header labels on row 1 so that Power BI can automatically turn
# the first row into a header when the table is imported.
#
# For most sheets the year labels (2015/2030/2040/2050) get prefixed with
# "Ano ". The "Potencia Incremental - SIN(MW)" sheet uses period ranges,
# so its labels get prefixed with "Intervalo " instead. The
# "Custo Total (bilhões de R$)" sheet only has a single year column (2015).

$wb = $excel.ActiveWorkbook

# Sheets whose header row uses plain years and gets the "Ano " prefix.
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet whose header row uses period ranges and gets the "Intervalo " prefix.
$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value = "Intervalo 2015"
$wsIntervalo.Range("C1").Value = "Intervalo 2015-2030"
$wsIntervalo.Range("D1").Value = "Intervalo 2031-2040"
$wsIntervalo.Range("E1").Value = "Intervalo 2041-2050"

# Sheet with a single year column.
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
